$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (old B..G shift to C..H; new I is added at the end)
$ws.Columns("B").Insert()

# Header row
$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = "Classifier"
$ws.Range("C1").Value = "Person"
$ws.Range("D1").Value = "Input #"
$ws.Range("E1").Value = "Input Song Name"
$ws.Range("F1").Value = "Input Song Artist"
$ws.Range("G1").Value = "Output Rec Name"
$ws.Range("H1").Value = "Output Rec Artist"
$ws.Range("I1").Value = "Assessment"

# Data rows (2-19)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Music Data"
$ws.Range("C2").Value = "Jon"
$ws.Range("D2").Value = 17924
$ws.Range("E2").Value = "Yonkers"
$ws.Range("F2").Value = "Tyler, The Creator"
$ws.Range("G2").Value = "Things to You"
$ws.Range("H2").Value = "Skeeter Davis, NRBQ"
$ws.Range("I2").Value = "Bad"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Music Data"
$ws.Range("C3").Value = "Jon"
$ws.Range("D3").Value = 30753
$ws.Range("E3").Value = "Crazy Little Thing Called Love"
$ws.Range("F3").Value = "Queen"
$ws.Range("G3").Value = "You"
$ws.Range("H3").Value = "George Duke"
$ws.Range("I3").Value = "Good"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Music Data"
$ws.Range("C4").Value = "Jon"
$ws.Range("D4").Value = 50106
$ws.Range("E4").Value = "Is This Love"
$ws.Range("F4").Value = "Bob Marley & The Wailers"
$ws.Range("G4").Value = "Drunk Off Your Love (feat. Sky Blu of LMFAO)"
$ws.Range("H4").Value = "''Shwayze', 'Cisco Adler'"
$ws.Range("I4").Value = "Good"

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Music Data"
$ws.Range("C5").Value = "Jon"
$ws.Range("D5").Value = 9625
$ws.Range("E5").Value = "Big Yellow Taxi"
$ws.Range("F5").Value = "Joni Mitchell"
$ws.Range("G5").Value = "Belief"
$ws.Range("H5").Value = "John Mayer"
$ws.Range("I5").Value = "Bad"

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Music Data"
$ws.Range("C6").Value = "Jon"
$ws.Range("D6").Value = 106893
$ws.Range("E6").Value = "Boom Boom Pow"
$ws.Range("F6").Value = "Black Eyed Peas"
$ws.Range("G6").Value = "La Barca Marina"
$ws.Range("H6").Value = "TIN TAN Y MARCELO"
$ws.Range("I6").Value = "Bad"

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Music Data"
$ws.Range("C7").Value = "Jon"
$ws.Range("D7").Value = 106889
$ws.Range("E7").Value = "Landslide"
$ws.Range("F7").Value = "''Stevie Nicks', 'Lindsey Buckingham'"
$ws.Range("G7").Value = "GCN Mushroom Bridge"
$ws.Range("H7").Value = "The Greatest Bits"
$ws.Range("I7").Value = "Bad"

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Music Data"
$ws.Range("C8").Value = "Jon"
$ws.Range("D8").Value = 112232
$ws.Range("E8").Value = "Do Nothing Till You Hear From Me"
$ws.Range("F8").Value = "Ella Fitzgerald"
$ws.Range("G8").Value = "Buona Sera"
$ws.Range("H8").Value = "Dean Martin"
$ws.Range("I8").Value = "Good"

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Music Data"
$ws.Range("C9").Value = "Jon"
$ws.Range("D9").Value = 19515
$ws.Range("E9").Value = "Old Town Road - Remix"
$ws.Range("F9").Value = "''Lil Nas X', 'Billy Ray Cyrus'"
$ws.Range("G9").Value = "Adiós Mi Amor"
$ws.Range("H9").Value = "Los Dareyes De La Sierra"
$ws.Range("I9").Value = "Bad"

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Music Data"
$ws.Range("C10").Value = "Jon"
$ws.Range("D10").Value = 19520
$ws.Range("E10").Value = "Dancing With A Stranger (with Normani)"
$ws.Range("F10").Value = "''Sam Smith', 'Normani'"
$ws.Range("G10").Value = "Dream Is Collapsing"
$ws.Range("H10").Value = "Hans Zimmer"
$ws.Range("I10").Value = "Bad"

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Music Data"
$ws.Range("C11").Value = "Jon"
$ws.Range("D11").Value = 30753
$ws.Range("E11").Value = "Crazy Little Thing Called Love"
$ws.Range("F11").Value = "Queen"
$ws.Range("G11").Value = "Dreams"
$ws.Range("H11").Value = "''The Rippingtons', 'Russ Freeman', 'David Benoit', 'David Koz', 'Gregg Karukas', 'Jimmy Johnson', 'Tony Morales', 'Steve Reid'"
$ws.Range("I11").Value = "Bad"

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Playlist"
$ws.Range("C12").Value = "Jon"
$ws.Range("D12").Value = 30753
$ws.Range("E12").Value = "Crazy Little Thing Called Love"
$ws.Range("F12").Value = "Queen"
$ws.Range("G12").Value = "Radioactive"
$ws.Range("H12").Value = "Imagine Dragons"
$ws.Range("I12").Value = "Bad"

$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Refined Playlist"
$ws.Range("C13").Value = "Jon"
$ws.Range("D13").Value = 30753
$ws.Range("E13").Value = "Crazy Little Thing Called Love"
$ws.Range("F13").Value = "Queen"
$ws.Range("G13").Value = "Ho Hey"
$ws.Range("H13").Value = "Lumineers"
$ws.Range("I13").Value = "Good"

$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Music Data"
$ws.Range("C14").Value = "Jon"
$ws.Range("D14").Value = 15851
$ws.Range("E14").Value = "Harder, Better, Faster, Stronger"
$ws.Range("F14").Value = "Daft Punk"
$ws.Range("G14").Value = "Istanbul (Not Constantinople)"
$ws.Range("H14").Value = "The Four Lads"
$ws.Range("I14").Value = "Bad"

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = "Playlist"
$ws.Range("C15").Value = "Jon"
$ws.Range("D15").Value = 15851
$ws.Range("E15").Value = "Harder, Better, Faster, Stronger"
$ws.Range("F15").Value = "Daft Punk"
$ws.Range("G15").Value = "Voyager"
$ws.Range("H15").Value = "Daft Punk"
$ws.Range("I15").Value = "Good"

$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Refined Playlist"
$ws.Range("C16").Value = "Jon"
$ws.Range("D16").Value = 15851
$ws.Range("E16").Value = "Harder, Better, Faster, Stronger"
$ws.Range("F16").Value = "Daft Punk"
$ws.Range("G16").Value = "Voyager"
$ws.Range("H16").Value = "Daft Punk"
$ws.Range("I16").Value = "Good"

$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "Music Data"
$ws.Range("C17").Value = "Jon"
$ws.Range("D17").Value = 9625
$ws.Range("E17").Value = "Big Yellow Taxi"
$ws.Range("F17").Value = "Joni Mitchell"
$ws.Range("G17").Value = "Belief"
$ws.Range("H17").Value = "John Mayer"
$ws.Range("I17").Value = "Bad"

$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "Playlist"
$ws.Range("C18").Value = "Jon"
$ws.Range("D18").Value = 9625
$ws.Range("E18").Value = "Big Yellow Taxi"
$ws.Range("F18").Value = "Joni Mitchell"
$ws.Range("G18").Value = "Radioactive"
$ws.Range("H18").Value = "Imagine Dragons"
$ws.Range("I18").Value = "Bad"

$ws.Range("A19").Value = 2
$ws.Range("B19").Value = "Refined Playlist"
$ws.Range("C19").Value = "Jon"
$ws.Range("D19").Value = 9625
$ws.Range("E19").Value = "Big Yellow Taxi"
$ws.Range("F19").Value = "Joni Mitchell"
$ws.Range("G19").Value = "Royals"
$ws.Range("H19").Value = "Lorde"
$ws.Range("I19").Value = "Bad"

